$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.289.52"
$ws.Range("E2").Value = "  +4.08%  "

$ws.Range("D3").Value = "3.485.34"
$ws.Range("E3").Value = "  +3.57%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'579.10"
$ws.Range("E5").Value = "  +3.35%  "

$ws.Range("D6").Value = "'162.31"
$ws.Range("E6").Value = "  +5.48%  "

$ws.Range("E7").Value = "  +14.73%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "3.487.78"
$ws.Range("E9").Value = "  +3.63%  "

$ws.Range("D10").Value = "'7.27"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").Value = "'0.126"
$ws.Range("E11").Value = "  +4.57%  "

$ws.Range("D12").Value = "'0.447"
$ws.Range("E12").Value = "  +4.18%  "

$ws.Range("D13").Value = "4.088.65"
$ws.Range("E13").Value = "  +3.63%  "

$ws.Range("D14").Value = "'0.134"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("E15").Value = "  +3.57%  "

$ws.Range("D16").Value = "'28.81"
$ws.Range("E16").Value = "  +7.34%  "

$ws.Range("D17").Value = "65.297.89"
$ws.Range("E17").Value = "  +3.97%  "

$ws.Range("D18").Value = "3.497.50"
$ws.Range("E18").Value = "  +6.53%  "

$ws.Range("E19").Value = "  +4.25%  "

$ws.Range("D20").Value = "'14.36"
$ws.Range("E20").Value = "  +2.91%  "

$ws.Range("D21").Value = "'384.06"
$ws.Range("E21").Value = "  +2.91%  "

$ws.Range("D22").Value = "'8.23"
$ws.Range("E22").Value = "  +3.36%  "

$ws.Range("D23").Value = "'0.551"
$ws.Range("E23").Value = "  +5.29%  "

$ws.Range("D24").Value = "'72.78"
$ws.Range("E24").Value = "  +2.77%  "

$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("D26").Value = "'0.0000120"
$ws.Range("E26").Value = "  +6.16%  "

$ws.Range("D27").Value = "'10.17"
$ws.Range("E27").Value = "  +7.55%  "

$ws.Range("E28").Value = "  +2.62%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("D30").Value = "'1.53"
$ws.Range("E30").Value = "  +14.21%  "

$ws.Range("D31").Value = "'6.26"
$ws.Range("E31").Value = "  +3.78%  "

$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  +4.81%  "

$ws.Range("D33").Value = "'23.64"
$ws.Range("E33").Value = "  +2.85%  "

$ws.Range("D34").Value = "'7.23"
$ws.Range("E34").Value = "  +7.98%  "

$ws.Range("D35").Value = "'1.59"
$ws.Range("E35").Value = "  +10.40%  "

$ws.Range("D36").Value = "'161.50"
$ws.Range("E36").Value = "  +1.50%  "

$ws.Range("D37").Value = "'1.93"
$ws.Range("E37").Value = "  +6.74%  "

$ws.Range("D38").Value = "3.048.86"
$ws.Range("E38").Value = "  +4.43%  "

$ws.Range("D39").Value = "'0.0776"
$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("D40").Value = "'27.09"
$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("D41").Value = "'4.58"
$ws.Range("E41").Value = "  +7.01%  "

$ws.Range("E42").Value = "  +2.15%  "

$ws.Range("D43").Value = "'6.53"
$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").Value = "'42.90"
$ws.Range("E44").Value = "  +4.10%  "

$ws.Range("D45").Value = "'0.780"
$ws.Range("E45").Value = "  +5.55%  "

$ws.Range("D46").Value = "'26.08"
$ws.Range("E46").Value = "  +13.87%  "

$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  +5.95%  "

$ws.Range("D48").Value = "'319.79"
$ws.Range("E48").Value = "  +13.20%  "

$ws.Range("D49").Value = "'6.77"
$ws.Range("E49").Value = "  +7.09%  "

$ws.Range("E50").Value = "  +7.91%  "

$ws.Range("D51").Value = "'2.22"
$ws.Range("E51").Value = "  +6.31%  "
